{"js": "// The author had written \"$495\" in the justification paragraph but the\n// correct (Government Staff) price quoted later in the document is\n// \"$595\" -- this fixes that missed price mention. Word also relocates\n// its auto-tracked \"_GoBack\" bookmark (last-edit position) from its old\n// spot to the point of this edit, so we move the bookmark to match.\n\nconst body = context.document.body;\n\n// 1) Correct the price mention: $495 -> $595\nconst priceMatches = body.search(\"$495\", { matchCase: true });\nawait context.sync();\n\nif (priceMatches.items.length === 0) {\n  throw new Error('Could not find \"$495\" to correct.');\n}\npriceMatches.items[0].insertText(\"$595\", \"Replace\");\nawait context.sync();\n\n// 2) Drop the \"_GoBack\" bookmark from its previous location (right before\n//    the trailing \"95\" of \"$1595\" later in the document), if present.\nconst oldBookmark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\nif (!oldBookmark.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 3) Re-insert \"_GoBack\" at the new edit point: right after \"$5\" and\n//    before \"95. \" in the sentence we just fixed.\nconst tailMatches = body.search(\"95. \", { matchCase: true });\nawait context.sync();\n\nif (tailMatches.items.length === 0) {\n  throw new Error('Could not find \"95. \" to anchor the bookmark.');\n}\nconst bookmarkPoint = tailMatches.items[0].getRange(\"Start\");\nbookmarkPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The author had written \"$495\" in the justification paragraph but the\n# correct (Government Staff) price quoted later in the document is\n# \"$595\" -- this fixes that missed price mention. Word also relocates\n# its auto-tracked \"_GoBack\" bookmark (last-edit position) from its old\n# spot to the point of this edit, so we move the bookmark to match.\n\n$d = $word.ActiveDocument\n\n# 1) Correct the price mention: $495 -> $595\n$find = $d.Content.Find\n$find.Text = \"`$495\"\n$find.Replacement.Text = \"`$595\"\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\nif (-not $found) {\n    throw 'Could not find \"$495\" to correct.'\n}\n\n# 2) Drop the \"_GoBack\" bookmark from its previous location (right before\n#    the trailing \"95\" of \"$1595\" later in the document).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 3) Re-insert \"_GoBack\" at the new edit point: right after \"$5\" and\n#    before \"95. \" in the sentence we just fixed.\n$r = $d.Content\n$found2 = $r.Find.Execute(\"95. \")\nif (-not $found2) {\n    throw 'Could not find \"95. \" to anchor the bookmark.'\n}\n$r.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $r)\n"}
